# Rename the "Query1" worksheet tab to "TotalsByProduct " (note trailing
# space) so it matches the underlying Power Query name ("TotalsByProduct").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Query1")
$ws.Name = "TotalsByProduct "
